function Set-TextValue {
    param($ws, $cellRef, $val)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws "D2" "67.350.50"
$ws.Range("E2").Value = "  -0.99%  "
Set-TextValue $ws "D3" "3.335.99"
$ws.Range("E3").Value = "  +2.20%  "
$ws.Range("E4").Value = "  -0.03%  "
Set-TextValue $ws "D5" "578.91"
$ws.Range("E5").Value = "  -0.77%  "
Set-TextValue $ws "D6" "184.23"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue $ws "D7" "0.605"
$ws.Range("E7").Value = "  +0.68%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue $ws "D8" "1.00"
$ws.Range("E8").Value = "  -0.02%  "
Set-TextValue $ws "D9" "0.129"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("E10").Value = "  +0.72%  "
Set-TextValue $ws "D11" "0.407"
Set-TextValue $ws "D12" "3.908.51"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("E13").Value = "  -0.87%  "
Set-TextValue $ws "D14" "27.34"
$ws.Range("E14").Value = "  -0.17%  "
Set-TextValue $ws "D15" "67.564.07"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("E16").Value = "  -0.25%  "
Set-TextValue $ws "D17" "3.326.62"
$ws.Range("E17").Value = "  +1.84%  "
Set-TextValue $ws "D18" "444.09"
$ws.Range("E18").Value = "  +6.35%  "
Set-TextValue $ws "D19" "13.57"
$ws.Range("E19").Value = "  +2.19%  "
Set-TextValue $ws "D20" "5.67"
$ws.Range("E20").Value = "  -0.79%  "
Set-TextValue $ws "D21" "7.70"
$ws.Range("E21").Value = "  +2.28%  "
Set-TextValue $ws "D22" "74.03"
$ws.Range("E22").Value = "  +4.07%  "
Set-TextValue $ws "D23" "0.999"
$ws.Range("E23").Value = "  -0.17%  "
Set-TextValue $ws "D24" "3.465.41"
$ws.Range("E24").Value = "  +1.84%  "
Set-TextValue $ws "D25" "0.512"
$ws.Range("E25").Value = "  +0.82%  "
Set-TextValue $ws "D26" "0.0000119"
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("E27").Value = "  +2.01%  "
Set-TextValue $ws "D28" "9.04"
$ws.Range("E28").Value = "  -3.64%  "
Set-TextValue $ws "D29" "1.01"
$ws.Range("E29").Value = "  +1.05%  "
Set-TextValue $ws "D30" "1.96"
$ws.Range("E30").Value = "  +0.83%  "
Set-TextValue $ws "D31" "22.86"
$ws.Range("E31").Value = "  +1.35%  "
Set-TextValue $ws "D32" "5.33"
$ws.Range("E32").Value = "  -2.25%  "
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("E36").Value = "  +4.25%  "
Set-TextValue $ws "D37" "161.25"
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("E38").Value = "  -1.94%  "
Set-TextValue $ws "D39" "27.08"
$ws.Range("E39").Value = "  +0.93%  "
Set-TextValue $ws "D40" "2.793.49"
$ws.Range("E40").Value = "  +6.06%  "
Set-TextValue $ws "D41" "0.791"
$ws.Range("E41").Value = "  -0.41%  "
Set-TextValue $ws "D42" "4.45"
$ws.Range("E42").Value = "  +0.15%  "
Set-TextValue $ws "D43" "6.23"
$ws.Range("E43").Value = "  -0.78%  "
Set-TextValue $ws "D44" "40.34"
$ws.Range("E44").Value = "  -0.78%  "
Set-TextValue $ws "D45" "0.0672"
$ws.Range("E45").Value = "  -0.26%  "
Set-TextValue $ws "D46" "24.60"
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("E47").Value = "  -2.44%  "
Set-TextValue $ws "D48" "325.45"
$ws.Range("E48").Value = "  -3.43%  "
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("E50").Value = "  +0.86%  "
Set-TextValue $ws "D51" "31.05"
$ws.Range("E51").Value = "  +1.80%  "
